$d = $word.ActiveDocument

# Color helper: Word COM uses BGR-ordered decimal color values.
# 0x7f6000 (RGB) -> BGR decimal = 0x00 * 65536 + 0x60 * 256 + 0x7f
$addColor = 24703   # 0x7f6000 as BGR-decimal (used for <add> tag runs)

# --- Change 1: "elle faict de grands traicts <del><add>se</add></del> tirants..."
# Split the "<del><add>" run so that "<add>" becomes its own run with color 7f6000.
$rng = $d.Content
$found = $rng.Find.Execute("<del><add>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find '<del><add>' run" }
$addPart = $d.Range($rng.Start + 5, $rng.End)
if ($addPart.Text -ne "<add>") { throw "Unexpected text for <add> split: [$($addPart.Text)]" }
$addPart.Font.Color = $addColor

# Split the "</add></del>" run so that "</add>" becomes its own run with color 7f6000.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("</add></del>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find '</add></del>' run" }
$addClosePart = $d.Range($rng2.Start, $rng2.Start + 6)
if ($addClosePart.Text -ne "</add>") { throw "Unexpected text for </add> split: [$($addClosePart.Text)]" }
$addClosePart.Font.Color = $addColor

# --- Change 2: "...d'u<add><exp>n</exp>/add> petit <tl>..."
# Split the "<add><exp>" run so that "<add>" becomes its own run with color 7f6000.
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("<add><exp>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Could not find '<add><exp>' run" }
$addPart2 = $d.Range($rng3.Start, $rng3.Start + 5)
if ($addPart2.Text -ne "<add>") { throw "Unexpected text for <add> split 2: [$($addPart2.Text)]" }
$addPart2.Font.Color = $addColor
# This particular <add> run originally sits inside a sz=14 (<exp>) run, but the
# canonical "<add>" tag formatting elsewhere in the document uses sz=18 (9pt).
$addPart2.Font.Size = 9

# Fix the typo "</exp>/add>" -> "</exp></add>" (no run split, just text correction).
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("</exp>/add>", $true, $false, $false, $false, $false, $true, 1, $false, "</exp></add>", 2)
if (-not $found4) { throw "Could not find/replace '</exp>/add>'" }
